$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Daily update: append the next day's row (row 49) following the existing
# pattern (date serial in column A, win counts in B:D).
$ws.Range("A49").Value = 45998
$ws.Range("B49").Value = 115
$ws.Range("C49").Value = 123
$ws.Range("D49").Value = 116

# Match the date-formatted style already used by the rest of column A.
$ws.Range("A49").NumberFormat = $ws.Range("A48").NumberFormat
